# login page code refactoring
# Re-orders the login test-data table (rows 2-5 cyclic shift), updates the
# "hyperlink preview" styling on A2, trims the set of mailto hyperlinks down
# to the two rows that still need them, and nudges a couple of cosmetic
# sheet-view settings (row height, selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-order the data rows (old r2 moves down to r5; r3,r4,r5 shift up) ---
# Capture the "before" values first so the writes below can't clobber data
# we still need to read.
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2

$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$c3 = $ws.Range("C3").Value2

$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2

$a5 = $ws.Range("A5").Value2
$b5 = $ws.Range("B5").Value2
$c5 = $ws.Range("C5").Value2

# Drop every mailto hyperlink; the two that survive are re-created below once
# the data has landed in its new rows (per-item .Delete() only mutates the
# in-memory collection view on this host, so clearing via the sheet-level
# collection is the reliable way to start from a clean slate).
$ws.Hyperlinks.Delete()

# New row 2 = old row 3
$ws.Range("A2").Value2 = $a3
$ws.Range("B2").Value2 = $b3
$ws.Range("C2").Value2 = $c3

# New row 3 = old row 4
$ws.Range("A3").Value2 = $a4
$ws.Range("B3").Value2 = $b4
$ws.Range("C3").Value2 = $c4

# New row 4 = old row 5
$ws.Range("A4").Value2 = $a5
$ws.Range("B4").Value2 = $b5
$ws.Range("C4").Value2 = $c5

# New row 5 = old row 2
$ws.Range("A5").Value2 = $a2
$ws.Range("B5").Value2 = $b2
$ws.Range("C5").Value2 = $c2

# --- 2. Re-create the two surviving hyperlinks (A2 and A5 now both resolve
#        to the "admin@yourstore.com" row) ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:admin@yourstore.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:admin@yourstore.com")

# --- 3. A2 swaps from the big 14pt "link row" look to the compact, default
#        11pt Hyperlink style (still centered, still boxed) ---
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Font.Size = 11

# --- 4. Cosmetic sheet-view nudges ---
$ws.Rows("1:5").RowHeight = 19
$ws.Range("C3").Select()
